$p = $ppt.ActivePresentation

# --- 1. Slide 3 ("Test Cases"): merge the split runs in the last bullet
#        ("Make sure last card " + "is King of Diamonds") into a single
#        run, and drop the stray endParaRPr left behind by the edit. ---
$s3 = $p.Slides.Item(3)
$shp = $s3.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Re-typing the whole body (through a distinct placeholder first so the
# host doesn't treat it as a no-op) normalizes run/endParaRPr bookkeeping
# for the edited paragraph. Re-apply the lvl=1 indent that two of the
# bullets use, since it is not otherwise preserved across the rewrite.
$tr.Text = "placeholder full text"
$full = "Make sure created deck has the correct number of cards.`rWrite if statement to see if the deck has the right number of cards`rIf not then tell go test handler that something is wrong.`rMake sure first card is Ace of Hearts`rMake sure last card is King of Diamonds"
$tr2 = $tf.TextRange
$tr2.Text = $full
$tr2.Paragraphs(2,1).IndentLevel = 2
$tr2.Paragraphs(3,1).IndentLevel = 2

# --- 2. Remove the trailing empty "Title/Content" slide (slide 4). ---
$p.Slides.Item(4).Delete()
